$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header renames (row 1) ---
$ws.Range("M1").Value = "FloodplainConnectivity_score"
$ws.Range("N1").Value = "Off-Channel/Side-Channels_score"

# --- Row 2 ---
$ws.Range("C2").Value = "Salmon Creek-Lower"
$ws.Range("M2").ClearContents()
$ws.Range("O2").ClearContents()

# --- Row 3 ---
$ws.Range("C3").Value = "Salmon Creek-Lower"
$ws.Range("M3").ClearContents()
$ws.Range("O3").ClearContents()

# --- Row 4 ---
$ws.Range("C4").Value = "Salmon Creek-Lower"
$ws.Range("M4").ClearContents()
$ws.Range("O4").ClearContents()

# --- Row 5 ---
$ws.Range("C5").Value = "Tonasket Creek DS"
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("O5").ClearContents()
